$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.311181333333333
$ws.Range("H2").Value = 6.933544
$ws.Range("I2").Value = 0.336172840858964
$ws.Range("J2").Value = 0.336172840858964
$ws.Range("M2").Value = 174.1282373333333
$ws.Range("N2").Value = 522.384712
$ws.Range("O2").Value = 0.985625830323027
$ws.Range("P2").Value = 0.985625830323027
$ws.Range("Q2").Value = 402.4419317310364
$ws.Range("R2").Value = 3621.977385579328
$ws.Range("S2").Value = 0.3313406354036672
$ws.Range("T2").Value = 0.3313406354036673
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.311181333333333
$ws.Range("H3").Value = 6.933544
$ws.Range("I3").Value = 0.336172840858964
$ws.Range("J3").Value = 0.336172840858964
$ws.Range("O3").Value = 0.003686901313133159
$ws.Range("P3").Value = 0.003686901313133159
$ws.Range("Q3").Value = 1.505402598948444
$ws.Range("R3").Value = 13.548623390536
$ws.Range("S3").Value = 0.001239436088402619
$ws.Range("T3").Value = 0.001239436088402619
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.311181333333333
$ws.Range("H4").Value = 6.933544
$ws.Range("I4").Value = 0.336172840858964
$ws.Range("J4").Value = 0.336172840858964
$ws.Range("M4").Value = 1.888095
$ws.Range("N4").Value = 5.664285
$ws.Range("O4").Value = 0.01068726836383999
$ws.Range("P4").Value = 0.01068726836383999
$ws.Range("Q4").Value = 4.36372991956
$ws.Range("R4").Value = 39.27356927604001
$ws.Range("S4").Value = 0.003592769366894223
$ws.Range("T4").Value = 0.003592769366894223
$ws.Range("I5").Value = 0.5963918049111226
$ws.Range("J5").Value = 0.5963918049111226
$ws.Range("M5").Value = 174.1282373333333
$ws.Range("N5").Value = 522.384712
$ws.Range("O5").Value = 0.985625830323027
$ws.Range("P5").Value = 0.985625830323027
$ws.Range("Q5").Value = 713.9573483203699
$ws.Range("R5").Value = 6425.616134883328
$ws.Range("S5").Value = 0.587819167913374
$ws.Range("T5").Value = 0.587819167913374
$ws.Range("I6").Value = 0.5963918049111226
$ws.Range("J6").Value = 0.5963918049111226
$ws.Range("O6").Value = 0.003686901313133159
$ws.Range("P6").Value = 0.003686901313133159
$ws.Range("S6").Value = 0.002198837728668673
$ws.Range("T6").Value = 0.002198837728668673
$ws.Range("I7").Value = 0.5963918049111226
$ws.Range("J7").Value = 0.5963918049111226
$ws.Range("M7").Value = 1.888095
$ws.Range("N7").Value = 5.664285
$ws.Range("O7").Value = 0.01068726836383999
$ws.Range("P7").Value = 0.01068726836383999
$ws.Range("Q7").Value = 7.741531874560001
$ws.Range("R7").Value = 69.67378687104001
$ws.Range("S7").Value = 0.006373799269080074
$ws.Range("T7").Value = 0.006373799269080074
$ws.Range("G8").Value = 0.4636166666666667
$ws.Range("H8").Value = 1.39085
$ws.Range("I8").Value = 0.06743535422991333
$ws.Range("J8").Value = 0.06743535422991333
$ws.Range("M8").Value = 174.1282373333333
$ws.Range("N8").Value = 522.384712
$ws.Range("O8").Value = 0.985625830323027
$ws.Range("P8").Value = 0.985625830323027
$ws.Range("Q8").Value = 80.72875296502224
$ws.Range("R8").Value = 726.5587766852001
$ws.Range("S8").Value = 0.06646602700598578
$ws.Range("T8").Value = 0.06646602700598578
$ws.Range("G9").Value = 0.4636166666666667
$ws.Range("H9").Value = 1.39085
$ws.Range("I9").Value = 0.06743535422991333
$ws.Range("J9").Value = 0.06743535422991333
$ws.Range("O9").Value = 0.003686901313133159
$ws.Range("P9").Value = 0.003686901313133159
$ws.Range("Q9").Value = 0.3019796520722223
$ws.Range("R9").Value = 2.71781686865
$ws.Range("S9").Value = 0.0002486274960618672
$ws.Range("T9").Value = 0.0002486274960618672
$ws.Range("G10").Value = 0.4636166666666667
$ws.Range("H10").Value = 1.39085
$ws.Range("I10").Value = 0.06743535422991333
$ws.Range("J10").Value = 0.06743535422991333
$ws.Range("M10").Value = 1.888095
$ws.Range("N10").Value = 5.664285
$ws.Range("O10").Value = 0.01068726836383999
$ws.Range("P10").Value = 0.01068726836383999
$ws.Range("Q10").Value = 0.8753523102500002
$ws.Range("R10").Value = 7.878170792250001
$ws.Range("S10").Value = 0.0007206997278656962
$ws.Range("T10").Value = 0.0007206997278656962
